$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = "四川"
$ws.Range("C19").Value = "资阳"
$ws.Range("B20").Value = "辽宁"
$ws.Range("C20").Value = "葫芦岛"
$ws.Range("B21").Value = "吉林"
$ws.Range("C21").Value = "长春"
$ws.Range("B22").Value = "四川"
$ws.Range("C22").Value = "广安"
$ws.Range("B30").Value = "河南"
$ws.Range("C30").Value = "洛阳"
$ws.Range("B31").Value = "山东"
$ws.Range("C31").Value = "滨州"
$ws.Range("B42").Value = "河北"
$ws.Range("C42").Value = "邯郸"
$ws.Range("B43").Value = "山东"
$ws.Range("C43").Value = "潍坊"
$ws.Range("B44").Value = "河南"
$ws.Range("C44").Value = "新乡"
$ws.Range("B45").Value = "河南"
$ws.Range("C45").Value = "安阳"
$ws.Range("B46").Value = "四川"
$ws.Range("C46").Value = "南充"
$ws.Range("B48").Value = "山西"
$ws.Range("C48").Value = "太原"
$ws.Range("B49").Value = "山东"
$ws.Range("C49").Value = "临沂"
$ws.Range("C53").Value = "德州"
$ws.Range("C54").Value = "济南"
$ws.Range("B58").Value = "安徽"
$ws.Range("C58").Value = "亳州"
$ws.Range("B59").Value = "山东"
$ws.Range("C59").Value = "日照"
$ws.Range("C60").Value = "枣庄"
$ws.Range("B63").Value = "安徽"
$ws.Range("C63").Value = "宿州"
$ws.Range("B64").Value = "陕西"
$ws.Range("C64").Value = "宝鸡"
$ws.Range("B68").Value = "江苏"
$ws.Range("C68").Value = "连云港"
$ws.Range("B69").Value = "河南"
$ws.Range("C69").Value = "南阳"
$ws.Range("B70").Value = "山西"
$ws.Range("C70").Value = "吕梁"
$ws.Range("B71").Value = "四川"
$ws.Range("C71").Value = "雅安"
$ws.Range("B72").Value = "河南"
$ws.Range("C72").Value = "郑州"
$ws.Range("B73").Value = "山东"
$ws.Range("C73").Value = "济宁"
$ws.Range("B78").Value = "河南"
$ws.Range("C78").Value = "许昌"
$ws.Range("B79").Value = "江苏"
$ws.Range("C79").Value = "宿迁"
$ws.Range("B80").Value = "河南"
$ws.Range("C80").Value = "开封"
$ws.Range("B81").Value = "四川"
$ws.Range("C81").Value = "泸州"
$ws.Range("B82").Value = "河南"
$ws.Range("C82").Value = "平顶山"
$ws.Range("B83").Value = "江苏"
$ws.Range("C83").Value = "徐州"
$ws.Range("B84").Value = "河南"
$ws.Range("C84").Value = "驻马店"
$ws.Range("B85").Value = "湖北"
$ws.Range("C85").Value = "武汉"
$ws.Range("B89").Value = "湖北"
$ws.Range("C89").Value = "孝感"
$ws.Range("B90").Value = "河南"
$ws.Range("C90").Value = "漯河"
$ws.Range("B91").Value = "山西"
$ws.Range("C91").Value = "朔州"
$ws.Range("B95").Value = "浙江"
$ws.Range("C95").Value = "湖州"
$ws.Range("B96").Value = "湖北"
$ws.Range("C96").Value = "随州"
$ws.Range("B97").Value = "江苏"
$ws.Range("C97").Value = "淮安"
$ws.Range("B98").Value = "四川"
$ws.Range("C98").Value = "遂宁"
$ws.Range("B99").Value = "湖北"
$ws.Range("C99").Value = "襄阳"
$ws.Range("B100").Value = "宁夏"
$ws.Range("C100").Value = "银川"
$ws.Range("B102").Value = "安徽"
$ws.Range("C102").Value = "淮北"
$ws.Range("B103").Value = "江苏"
$ws.Range("C103").Value = "常州"
$ws.Range("B104").Value = "湖南"
$ws.Range("C104").Value = "益阳"
$ws.Range("B108").Value = "江苏"
$ws.Range("C108").Value = "盐城"
$ws.Range("B109").Value = "安徽"
$ws.Range("C109").Value = "六安"
$ws.Range("C110").Value = "蚌埠"
$ws.Range("B111").Value = "湖南"
$ws.Range("C111").Value = "岳阳"
$ws.Range("B112").Value = "青海"
$ws.Range("C112").Value = "西宁"
$ws.Range("B116").Value = "广东"
$ws.Range("C116").Value = "广州"
$ws.Range("B117").Value = "新疆"
$ws.Range("C117").Value = "乌鲁木齐"
$ws.Range("B118").Value = "贵州"
$ws.Range("C118").Value = "贵阳"
$ws.Range("B119").Value = "山西"
$ws.Range("C119").Value = "大同"
$ws.Range("B120").Value = "江西"
$ws.Range("C120").Value = "新余"
$ws.Range("B121").Value = "湖北"
$ws.Range("C121").Value = "荆门"
$ws.Range("B122").Value = "安徽"
$ws.Range("C122").Value = "铜陵"
$ws.Range("B123").Value = "江苏"
$ws.Range("C123").Value = "扬州"
$ws.Range("B124").Value = "湖北"
$ws.Range("C124").Value = "鄂州"
$ws.Range("B125").Value = "湖南"
$ws.Range("C125").Value = "常德"
$ws.Range("B126").Value = "江苏"
$ws.Range("C126").Value = "无锡"
$ws.Range("B127").Value = "湖北"
$ws.Range("C127").Value = "荆州"
$ws.Range("B128").Value = "江苏"
$ws.Range("C128").Value = "南通"
$ws.Range("B129").Value = "湖北"
$ws.Range("C129").Value = "咸宁"
$ws.Range("B130").Value = "湖南"
$ws.Range("C130").Value = "长沙"
$ws.Range("B133").Value = "辽宁"
$ws.Range("C133").Value = "大连"
$ws.Range("B134").Value = "湖北"
$ws.Range("C134").Value = "黄冈"
$ws.Range("B136").Value = "湖北"
$ws.Range("C136").Value = "黄石"
$ws.Range("B137").Value = "浙江"
$ws.Range("C137").Value = "绍兴"
$ws.Range("B138").Value = "福建"
$ws.Range("C138").Value = "厦门"
$ws.Range("B146").Value = "江苏"
$ws.Range("C146").Value = "南京"
$ws.Range("B147").Value = "湖南"
$ws.Range("C147").Value = "湘潭"
$ws.Range("B159").Value = "浙江"
$ws.Range("C159").Value = "台州"
$ws.Range("B160").Value = "广西"
$ws.Range("C160").Value = "南宁"
